$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("E3").Value = "-"

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("E4").Value = "[-, -, -, 'MCT-3A-Eletrohidráulica']"

# Row 6
$ws.Range("B6").Value = "-"
$ws.Range("E6").Value = "[-, -, -, 'MCT-3A-Eletrohidráulica']"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("E7").Value = "[-, -, -, 'MCT-3A-Eletrohidráulica']"

# Row 8
$ws.Range("E8").Value = "[-, -, -, 'MCT-3A-Eletrohidráulica']"

# Row 11
$ws.Range("B11").Value = "-"
$ws.Range("C11").Value = "-"

# Row 12
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "-"

# Row 14
$ws.Range("C14").Value = "-"

# Row 15
$ws.Range("C15").Value = "-"
